# Regenerate sval data to filter save games: replace the TB/d2S/K/IP
# columns with their recomputed values and refresh the sum column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (TB), C (d2S), D (K), E (IP); column G (sum)
# is recomputed as B+C+D+E for every data row (rows 2-18).
$data = @(
    @{ Row = 2;  B = 3.272327238179451;  C = 1.626987699542094;  D = 0.1496068669990043;  E = 0.5333859586016987 },
    @{ Row = 3;  B = 3.272327238179451;  C = 1.626987699542094;  D = 18.71679738969934;   E = 0.5333859586016987 },
    @{ Row = 4;  B = 3.272327238179451;  C = 1.626987699542094;  D = 0.7210945179870265;  E = 0.5333859586016987 },
    @{ Row = 5;  B = 0.6545652718822623; C = 1.626987699542094;  D = 0.1496068669990043;  E = 0.5333859586016987 },
    @{ Row = 6;  B = 0.2881169905109251; C = 0.3048912486333797; D = 3.223369029078222;   E = 0.5333859586016987 },
    @{ Row = 7;  B = 3.272327238179451;  C = 1.626987699542094;  D = 0.7210945179870265;  E = 0.5333859586016987 },
    @{ Row = 8;  B = 0.1169995834814548; C = 0.04103571897497393;D = 3.223369029078222;   E = 0.5333859586016987 },
    @{ Row = 9;  B = 0.1169995834814548; C = 0.3048912486333797; D = 0.7210945179870265;  E = 0.5333859586016987 },
    @{ Row = 10; B = 3.272327238179451;  C = 1.626987699542094;  D = 0.7210945179870265;  E = 0.5333859586016987 },
    @{ Row = 11; B = 3.272327238179451;  C = 1.626987699542094;  D = 0.7210945179870265;  E = 0.5333859586016987 },
    @{ Row = 12; B = 3.272327238179451;  C = 1.626987699542094;  D = 0.7210945179870265;  E = 0.5333859586016987 },
    @{ Row = 13; B = 3.272327238179451;  C = 1.626987699542094;  D = 0.7210945179870265;  E = 0.5333859586016987 },
    @{ Row = 14; B = 3.272327238179451;  C = 1.626987699542094;  D = 3.223369029078222;   E = 0.5333859586016987 },
    @{ Row = 15; B = 3.272327238179451;  C = 1.626987699542094;  D = 0.7210945179870265;  E = 13.86384647080068 },
    @{ Row = 16; B = 1.445647641019636;  C = 9.983522426115931;  D = 0.1496068669990043;  E = 13.86384647080068 },
    @{ Row = 17; B = 0.6545652718822623; C = 1.626987699542094;  D = 0.7210945179870265;  E = 0.5333859586016987 },
    @{ Row = 18; B = 1.445647641019636;  C = 1.626987699542094;  D = 3.223369029078222;   E = 0.5333859586016987 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Range("B$r").Value = $entry.B
    $ws.Range("C$r").Value = $entry.C
    $ws.Range("D$r").Value = $entry.D
    $ws.Range("E$r").Value = $entry.E
    $ws.Range("G$r").Value = $entry.B + $entry.C + $entry.D + $entry.E
}
